# This script re-runs the projection model outputs on the "Comparison" sheet
# with a new random seed for termination dates, and adds prorated compensation
# for employees who terminated mid-year. This changes headcount/eligibility/
# participation counts (cols C:G) and all downstream dollar totals and ratios
# (cols H:R) for every scenario/year row (2-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Baseline / Year 1
$ws.Range("D2").Value = 101
$ws.Range("E2").Value = 85
$ws.Range("F2").Value = 0.8415841584158416
$ws.Range("G2").Value = 0.8333333333333334
$ws.Range("H2").Value = 0.09943492116856165
$ws.Range("I2").Value = 0.0828624343071347
$ws.Range("J2").Value = 453942.9050599152
$ws.Range("K2").Value = 164357.4527309576
$ws.Range("M2").Value = 164357.4527309576
$ws.Range("N2").Value = 618300.3577908728
$ws.Range("O2").Value = 10015935.0988
$ws.Range("P2").Value = 9608194.168699998
$ws.Range("Q2").Value = 0.01640959641907515
$ws.Range("R2").Value = 0.01710596703659199

# Row 3: Baseline / Year 2
$ws.Range("D3").Value = 103
$ws.Range("F3").Value = 0.8446601941747572
$ws.Range("H3").Value = 0.09776072146367178
$ws.Range("I3").Value = 0.08257458997416935
$ws.Range("J3").Value = 475461.5553898957
$ws.Range("K3").Value = 172518.3579019778
$ws.Range("M3").Value = 172518.3579019778
$ws.Range("N3").Value = 647979.9132918735
$ws.Range("O3").Value = 10570506.655464
$ws.Range("P3").Value = 10163133.497461
$ws.Range("Q3").Value = 0.01632072742821664
$ws.Range("R3").Value = 0.01697491801569636

# Row 4: Baseline / Year 3
$ws.Range("E4").Value = 88
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.8461538461538461
$ws.Range("H4").Value = 0.09666630872252636
$ws.Range("I4").Value = 0.08179456891906077
$ws.Range("J4").Value = 499772.3434468232
$ws.Range("K4").Value = 177587.3793366524
$ws.Range("M4").Value = 177587.3793366524
$ws.Range("N4").Value = 677359.7227834756
$ws.Range("O4").Value = 10935308.90212792
$ws.Range("P4").Value = 10526964.54938483
$ws.Range("Q4").Value = 0.01623981370129338
$ws.Range("R4").Value = 0.01686976131662097

# Row 5: Baseline / Year 4
$ws.Range("E5").Value = 90
$ws.Range("F5").Value = 0.8571428571428571
$ws.Range("G5").Value = 0.8571428571428571
$ws.Range("H5").Value = 0.09525529745848933
$ws.Range("I5").Value = 0.0816473978215623
$ws.Range("J5").Value = 519331.4335515244
$ws.Range("K5").Value = 184792.9606174003
$ws.Range("M5").Value = 184792.9606174003
$ws.Range("N5").Value = 704124.3941689247
$ws.Range("O5").Value = 11301297.20929176
$ws.Range("P5").Value = 10890602.52596637
$ws.Range("Q5").Value = 0.0163514822409472
$ws.Range("R5").Value = 0.01696811174375338

# Row 6: Baseline / Year 5
$ws.Range("H6").Value = 0.0962254754163797
$ws.Range("I6").Value = 0.0798853003456737
$ws.Range("J6").Value = 530251.3219092456
$ws.Range("K6").Value = 187191.7221115101
$ws.Range("M6").Value = 187191.7221115101
$ws.Range("N6").Value = 717443.0440207556
$ws.Range("O6").Value = 11764042.88127051
$ws.Range("P6").Value = 11349577.35744537
$ws.Range("Q6").Value = 0.01591219311258524
$ws.Range("R6").Value = 0.01649327690503925

# Row 7: AIP_New_Hires / Year 1
$ws.Range("D7").Value = 101
$ws.Range("E7").Value = 87
$ws.Range("F7").Value = 0.8613861386138614
$ws.Range("G7").Value = 0.8529411764705882
$ws.Range("H7").Value = 0.09847934560753153
$ws.Range("I7").Value = 0.0840951281162279
$ws.Range("J7").Value = 461526.6616369847
$ws.Range("K7").Value = 168149.3310194924
$ws.Range("M7").Value = 168149.3310194924
$ws.Range("N7").Value = 629675.9926564771
$ws.Range("O7").Value = 10098564.6888
$ws.Range("P7").Value = 9690823.758699998
$ws.Range("Q7").Value = 0.01665081486342128
$ws.Range("R7").Value = 0.01735139707483951

# Row 8: AIP_New_Hires / Year 2
$ws.Range("H8").Value = 0.09664744730659888
$ws.Range("I8").Value = 0.08351090107075052
$ws.Range("J8").Value = 480378.4435696948
$ws.Range("K8").Value = 174976.8019918774
$ws.Range("M8").Value = 174976.8019918774
$ws.Range("N8").Value = 655355.2455615721
$ws.Range("O8").Value = 10549636.406764
$ws.Range("P8").Value = 10142263.248761
$ws.Range("Q8").Value = 0.0165860504803454
$ws.Range("R8").Value = 0.01725224416880058

# Row 9: AIP_New_Hires / Year 3
$ws.Range("D9").Value = 104
$ws.Range("E9").Value = 90
$ws.Range("F9").Value = 0.8653846153846154
$ws.Range("G9").Value = 0.8653846153846154
$ws.Range("H9").Value = 0.09672223673719484
$ws.Range("I9").Value = 0.08370193563795707
$ws.Range("J9").Value = 511686.5069044705
$ws.Range("K9").Value = 183544.4610654761
$ws.Range("M9").Value = 183544.4610654761
$ws.Range("N9").Value = 695230.9679699468
$ws.Range("O9").Value = 10895668.70626692
$ws.Range("P9").Value = 10487324.35352383
$ws.Range("Q9").Value = 0.01684563527155573
$ws.Range("R9").Value = 0.01750155281540459

# Row 10: AIP_New_Hires / Year 4
$ws.Range("H10").Value = 0.09467907571674584
$ws.Range("I10").Value = 0.08205519895451308
$ws.Range("J10").Value = 524337.9468682207
$ws.Range("K10").Value = 187296.2172757485
$ws.Range("M10").Value = 187296.2172757485
$ws.Range("N10").Value = 711634.1641439691
$ws.Range("O10").Value = 11320413.25095493
$ws.Range("P10").Value = 10909718.56762955
$ws.Range("Q10").Value = 0.01654499823669858
$ws.Range("R10").Value = 0.01716783215943617

# Row 11: AIP_New_Hires / Year 5
$ws.Range("D11").Value = 105
$ws.Range("E11").Value = 91
$ws.Range("F11").Value = 0.8666666666666667
$ws.Range("G11").Value = 0.8584905660377359
$ws.Range("H11").Value = 0.09657166851730856
$ws.Range("I11").Value = 0.08300020599127433
$ws.Range("J11").Value = 550804.6883894347
$ws.Range("K11").Value = 197468.4053516046
$ws.Range("M11").Value = 197468.4053516046
$ws.Range("N11").Value = 748273.0937410393
$ws.Range("O11").Value = 11614789.51518358
$ws.Range("P11").Value = 11200323.99135843
$ws.Range("Q11").Value = 0.01700146223859344
$ws.Range("R11").Value = 0.01763059760628002

# Row 12: AIP_All_Eligible / Year 1
$ws.Range("D12").Value = 101
$ws.Range("E12").Value = 101
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9901960784313726
$ws.Range("H12").Value = 0.09250231417550324
$ws.Range("I12").Value = 0.09159542874241007
$ws.Range("J12").Value = 538620.5426134155
$ws.Range("K12").Value = 206696.2715077078
$ws.Range("M12").Value = 206696.2715077078
$ws.Range("N12").Value = 745316.8141211235
$ws.Range("O12").Value = 10110406.7188
$ws.Range("P12").Value = 9702665.7887
$ws.Range("Q12").Value = 0.02044391261959444
$ws.Range("R12").Value = 0.02130303939237319

# Row 13: AIP_All_Eligible / Year 2
$ws.Range("D13").Value = 102
$ws.Range("E13").Value = 102
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9902912621359223
$ws.Range("H13").Value = 0.09875808338513453
$ws.Range("I13").Value = 0.09779926704158953
$ws.Range("J13").Value = 638784.0023400235
$ws.Range("K13").Value = 254179.5813770418
$ws.Range("M13").Value = 254179.5813770418
$ws.Range("N13").Value = 892963.5837170655
$ws.Range("O13").Value = 10480889.176664
$ws.Range("P13").Value = 10073516.018661
$ws.Range("Q13").Value = 0.02425171920937585
$ws.Range("R13").Value = 0.02523245914397504

# Row 14: AIP_All_Eligible / Year 3
$ws.Range("D14").Value = 104
$ws.Range("E14").Value = 104
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0.1018639034285131
$ws.Range("I14").Value = 0.1018639034285131
$ws.Range("J14").Value = 719445.6476580924
$ws.Range("K14").Value = 287424.0314422871
$ws.Range("M14").Value = 287424.0314422871
$ws.Range("N14").Value = 1006869.67910038
$ws.Range("O14").Value = 10862870.64886392
$ws.Range("P14").Value = 10454526.29612083
$ws.Range("Q14").Value = 0.0264593071880449
$ws.Range("R14").Value = 0.02749278382406827

# Row 15: AIP_All_Eligible / Year 4
$ws.Range("E15").Value = 105
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.1034799091691537
$ws.Range("I15").Value = 0.1034799091691537
$ws.Range("J15").Value = 770434.6693487333
$ws.Range("K15").Value = 310344.5785160048
$ws.Range("M15").Value = 310344.5785160048
$ws.Range("N15").Value = 1080779.247864738
$ws.Range("O15").Value = 11340380.44192984
$ws.Range("P15").Value = 10929685.75860446
$ws.Range("Q15").Value = 0.02736632867875746
$ws.Range("R15").Value = 0.02839464787646656

# Row 16: AIP_All_Eligible / Year 5
$ws.Range("D16").Value = 106
$ws.Range("E16").Value = 106
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0.1038456560215263
$ws.Range("I16").Value = 0.1038456560215263
$ws.Range("J16").Value = 810410.2254774929
$ws.Range("K16").Value = 327271.1738956337
$ws.Range("M16").Value = 327271.1738956337
$ws.Range("N16").Value = 1137681.399373127
$ws.Range("O16").Value = 11677572.73988773
$ws.Range("P16").Value = 11263107.21606259
$ws.Range("Q16").Value = 0.02802561638325363
$ws.Range("R16").Value = 0.0290569172092142
